$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 (current last row) reverts from the "latest row" date-only style to the normal datetime style
$ws.Range("A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 44 with the "latest row" date-only style
$ws.Range("A44").Value = 45629
$ws.Range("A44").NumberFormat = "YYYY-MM-DD"
$ws.Range("B44").Value = 115
$ws.Range("C44").Value = 98
$ws.Range("D44").Value = 105
